$d = $word.ActiveDocument
$paras = $d.Paragraphs

# ---------------------------------------------------------------------------
# Change 1: "Equipes" -> "Équipes" in the "Equipes : Toutes les ..." line,
# and drop the surrounding spell-check proofErr markers around that word.
# ---------------------------------------------------------------------------
$pEquipes = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "Equipes : Toutes les*") {
        $pEquipes = $paras.Item($i)
        break
    }
}

if ($pEquipes -ne $null) {
    $xmlEquipes = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1918BBA9" w14:textId="77777777" w:rsidR="00EF1450" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>Équipes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> : Toutes les équipes déclarées </w:t></w:r></w:p>'
    [void]$pEquipes.Range.InsertXML($xmlEquipes)
} else {
    Write-Output "WARNING: 'Equipes' paragraph not found"
}

# ---------------------------------------------------------------------------
# Change 2: in the "L'idée Algorithmique : ..." paragraph, split the run so
# "QR code" becomes "QR codes" (three runs), and replace the following empty
# paragraph with: an empty paragraph, a new underlined "Exemple de licence
# amateur :" paragraph, and a trailing empty underlined paragraph.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$pIdee = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "L?idée Algorithmique :*") {
        $pIdee = $i
        break
    }
}

if ($pIdee -ne $null) {
    $start = $paras.Item($pIdee).Range.Start
    $endp = $paras.Item($pIdee + 1).Range.End
    $r = $d.Range($start, $endp)

    $xmlIdee = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7727979D" w14:textId="77777777" w:rsidR="00EF1450" w:rsidRDefault="00000000"><w:pPr><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>L’idée Algorithmique : à partir des licences (données) joueurs, entraîneurs générer des QR code</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> (résultats) et pouvoir les vérifier via l’appareil photo (une vérification entre le message retourné par le QR code et une base de données de toutes les licences sera effectuée).</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:u w:val="single"/></w:rPr><w:t>Exemple de licence amateur :</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:u w:val="single"/></w:rPr></w:pPr></w:p>'

    [void]$r.InsertXML($xmlIdee)
} else {
    Write-Output "WARNING: 'L'idée Algorithmique' paragraph not found"
}

Write-Output "DONE"
